# Update the 25 two-digit x two-digit multiplication problems/answers
# in the answers table to new problem/answer pairs.
$d = $word.ActiveDocument

$d.Content.Find.Execute("39×48=1872", $true, $false, $false, $false, $false, $true, 1, $false, "36×53=1908", 2) | Out-Null
$d.Content.Find.Execute("99×66=6534", $true, $false, $false, $false, $false, $true, 1, $false, "73×51=3723", 2) | Out-Null
$d.Content.Find.Execute("44×56=2464", $true, $false, $false, $false, $false, $true, 1, $false, "73×93=6789", 2) | Out-Null
$d.Content.Find.Execute("87×82=7134", $true, $false, $false, $false, $false, $true, 1, $false, "27×81=2187", 2) | Out-Null
$d.Content.Find.Execute("28×15=420", $true, $false, $false, $false, $false, $true, 1, $false, "41×21=861", 2) | Out-Null
$d.Content.Find.Execute("73×96=7008", $true, $false, $false, $false, $false, $true, 1, $false, "87×79=6873", 2) | Out-Null
$d.Content.Find.Execute("41×37=1517", $true, $false, $false, $false, $false, $true, 1, $false, "26×80=2080", 2) | Out-Null
$d.Content.Find.Execute("94×48=4512", $true, $false, $false, $false, $false, $true, 1, $false, "23×75=1725", 2) | Out-Null
$d.Content.Find.Execute("74×49=3626", $true, $false, $false, $false, $false, $true, 1, $false, "97×94=9118", 2) | Out-Null
$d.Content.Find.Execute("29×86=2494", $true, $false, $false, $false, $false, $true, 1, $false, "27×47=1269", 2) | Out-Null
$d.Content.Find.Execute("17×58=986", $true, $false, $false, $false, $false, $true, 1, $false, "97×11=1067", 2) | Out-Null
$d.Content.Find.Execute("92×70=6440", $true, $false, $false, $false, $false, $true, 1, $false, "31×88=2728", 2) | Out-Null
$d.Content.Find.Execute("76×98=7448", $true, $false, $false, $false, $false, $true, 1, $false, "66×86=5676", 2) | Out-Null
$d.Content.Find.Execute("85×90=7650", $true, $false, $false, $false, $false, $true, 1, $false, "21×25=525", 2) | Out-Null
$d.Content.Find.Execute("32×92=2944", $true, $false, $false, $false, $false, $true, 1, $false, "46×36=1656", 2) | Out-Null
$d.Content.Find.Execute("12×99=1188", $true, $false, $false, $false, $false, $true, 1, $false, "92×12=1104", 2) | Out-Null
$d.Content.Find.Execute("13×70=910", $true, $false, $false, $false, $false, $true, 1, $false, "43×12=516", 2) | Out-Null
$d.Content.Find.Execute("94×96=9024", $true, $false, $false, $false, $false, $true, 1, $false, "75×44=3300", 2) | Out-Null
$d.Content.Find.Execute("26×45=1170", $true, $false, $false, $false, $false, $true, 1, $false, "38×15=570", 2) | Out-Null
$d.Content.Find.Execute("72×58=4176", $true, $false, $false, $false, $false, $true, 1, $false, "39×43=1677", 2) | Out-Null
$d.Content.Find.Execute("97×69=6693", $true, $false, $false, $false, $false, $true, 1, $false, "91×56=5096", 2) | Out-Null
$d.Content.Find.Execute("55×30=1650", $true, $false, $false, $false, $false, $true, 1, $false, "69×26=1794", 2) | Out-Null
$d.Content.Find.Execute("30×25=750", $true, $false, $false, $false, $false, $true, 1, $false, "43×66=2838", 2) | Out-Null
$d.Content.Find.Execute("66×44=2904", $true, $false, $false, $false, $false, $true, 1, $false, "45×56=2520", 2) | Out-Null
$d.Content.Find.Execute("15×59=885", $true, $false, $false, $false, $false, $true, 1, $false, "93×93=8649", 2) | Out-Null
